# Split two long run texts into multiple <w:t> runs separated by manual
# line breaks (<w:br/>), matching the authored edit.

$d = $word.ActiveDocument

# 1) "Critério:" paragraph text — insert a line break between the two
#    sentences (after "...na engenharia." and before "Para a formação...").
$d.Content.Find.Execute(
    "despertar interesse na engenharia.Para a formação",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "despertar interesse na engenharia.^lPara a formação", 2) | Out-Null

# 2) "Norma de recuperação:" paragraph text — insert a line break before
#    each "- " bullet-like item that follows a period (".- " -> ".<br>- ").
$d.Content.Find.Execute(
    ".- ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    ".^l- ", 2) | Out-Null
